# Finished making all sections populate into excel worksheet.
# The "Education" sheet had duplicate Training/Posting/Medal columns
# (O:AK) that are now properly split out into their own
# "Training and Posting info" sheet, so that duplicated data is
# cleared out of "Education" here.

$wb = $excel.ActiveWorkbook
$wsEducation = $wb.Worksheets.Item("Education")
$wsTraining  = $wb.Worksheets.Item("Training and Posting info")

# Clear out the duplicated training/posting/medal columns (O:AK) from
# the Education sheet - this removes the header labels in row 1 and
# all the data values in rows 2:15, leaving behind only the empty,
# still-formatted date cells (Q, W, X, AF, AG) that carry a date
# number format.
$wsEducation.Range("O1:AK15").ClearContents()

# Restore/record the view state on each sheet: leave the "Training and
# Posting info" sheet scrolled/selected at F4 (no longer the active
# tab)...
$wsTraining.Activate()
$wsTraining.Range("F4").Select()

# ...and make "Education" the active sheet again, with the former
# O:AK selection left selected (artifact of the cleanup) and the view
# scrolled down a bit.
$wsEducation.Activate()
$wsEducation.Range("O1:AK1048576").Select()
